$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (96) with the latest profit figure for 11/21/2025.
# The date is stored as literal text (matching the existing rows), so we
# temporarily force a text number format to stop Excel from converting the
# "MM/DD/YYYY" string into a date serial number, then clear the formatting
# afterwards so the new cells end up with the same (default) style as the
# rest of the data rows.
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = "11/21/2025"
$ws.Range("A96").ClearFormats()

$ws.Range("B96").Value = 6846.43
